$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A101").Value = 88.45342611741277
$ws.Range("B2:B101").Value = 41399.9172001656
$ws.Range("C2:C101").Value = 569.0215447350218
$ws.Range("D2:D101").Value = 2467.006475771309
$ws.Range("F2:F101").Value = 45055.898259560185
